# Applies the "Add files via upload" commit:
#   - Rename sheet "MGBC-MBL-01" -> "MGBC_BOM" (defined name / filter
#     database reference follows automatically since it targets the sheet).
#   - Update the Mouser cart-link URL stored in cell G1 of that sheet
#     (AccessID changed from 2362d7198d to 09d74b6d14).
#   - Bump the sheet's zoom level from 85% to 115%.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MGBC-MBL-01")

# Rename the worksheet; the workbook-level _FilterDatabase defined name
# tracks the sheet by reference, so it updates automatically.
$ws.Name = "MGBC_BOM"

# Update the "Cart link:" URL in G1.
$ws.Range("G1").Value = "https://www.mouser.com/ProjectManager/ProjectDetail.aspx?AccessID=09d74b6d14"

# Make sure this sheet is active, then set the zoom to 115%.
$ws.Activate()
$excel.ActiveWindow.Zoom = 115
